$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, shifting existing rows 156:241 down to 157:242
$ws.Rows("156:156").Insert()

# Populate the newly inserted row 156 with the new record's data
$ws.Cells.Item(156, 1).Value = 3
$ws.Cells.Item(156, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 45029
$ws.Cells.Item(156, 5).Value = 5
$ws.Cells.Item(156, 6).Value = 100112026
$ws.Cells.Item(156, 7).Value = "Haba"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 40
$ws.Cells.Item(156, 11).Value = 20000
$ws.Cells.Item(156, 12).Value = 20000
$ws.Cells.Item(156, 13).Value = 20000
$ws.Cells.Item(156, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(156, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(156, 16).Value = 800
$ws.Cells.Item(156, 17).Value = 25
$ws.Cells.Item(156, 18).Value = "Hortaliza"
